$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2 (Germany / "0 o 1" note) was removed from the sheet ---
# Clear out the old data row entirely; D2/E2 keep their date formatting
# but become blank, same as the rest of the row.
$ws.Range("A2:G2").ClearContents()

# --- US row 8: the survey start date was corrected ---
$ws.Range("D8").Value = 43909

# --- New row: Sweden, added at the bottom of the table ---
$ws.Range("A18").Value = "sweden"
$ws.Range("B18").Value = 0.85
$ws.Range("C18").Value = 0.1

# Reuse the date-formatted style already used by the other FechaInicial /
# FechaFinal columns (copy format only, then set the values).
$ws.Range("D3").Copy()
$ws.Range("D18:E18").PasteSpecial(-4122)
$ws.Application.CutCopyMode = $false
$ws.Range("D18").Value = 43867
$ws.Range("E18").Value = 43960

$ws.Range("G18").Value = "no"

# Source link for the new row.
$ws.Range("H18").Value = "https://www.google.com/covid19/mobility/"
$ws.Hyperlinks.Add($ws.Range("H18"), "https://www.google.com/covid19/mobility/")
$ws.Range("H18").Style = "Hyperlink"

# --- Update the saved selection / scroll position ---
$ws.Range("H8:H9").Select()
